# Atualização de bases das ligas, do dia: 27-04-2024 às 11:27
#
# The underlying source data had a handful of rows whose content (every
# column except the leading "id" column A) was swapped with an adjacent
# row. This applies that swap for the affected row pairs by exchanging
# the B:AB range contents between each pair of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($worksheet, $rowA, $rowB)

    $rangeA = $worksheet.Range("B$rowA" + ":AB$rowA")
    $rangeB = $worksheet.Range("B$rowB" + ":AB$rowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

Swap-RowRange $ws 104 105
Swap-RowRange $ws 107 108
Swap-RowRange $ws 144 145
Swap-RowRange $ws 148 150
Swap-RowRange $ws 211 212
Swap-RowRange $ws 214 215
